# Edit: insert a new daily price record at row 95 for
# "Femacal de La Calera - Arándano (blue)" (Fruta / hortaliza, semanal).
# Inserting the row shifts the existing row 95 (and everything below it)
# down by one; the freshly inserted row starts blank, so every column of
# the new record is written explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 95 (and all rows below) down by one to make room for the new record.
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the new record's data.
$ws.Cells.Item(95, 1).Value = 3                          # A95 - Mercado ID
$ws.Cells.Item(95, 2).Value = "Femacal de La Calera"     # B95 - Mercado
$ws.Cells.Item(95, 3).Value = "Coquimbo"                 # C95 - Región
$ws.Cells.Item(95, 4).Value = 44586                      # D95 - Fecha
$ws.Cells.Item(95, 5).Value = 5                           # E95 - Codreg
$ws.Cells.Item(95, 6).Value = "Fruta"                    # F95 - Tipo
$ws.Cells.Item(95, 7).Value = 100101                      # G95 - Producto ID
$ws.Cells.Item(95, 8).Value = "Berries"                  # H95 - Producto
$ws.Cells.Item(95, 9).Value = 100101001                   # I95 - Categoría ID
$ws.Cells.Item(95, 10).Value = "Arándano (blue)"         # J95 - Categoría
$ws.Cells.Item(95, 11).Value = "Sin especificar"         # K95 - Variedad
$ws.Cells.Item(95, 12).Value = "Primera"                 # L95 - Calidad
$ws.Cells.Item(95, 13).Value = 185                        # M95 - Volumen
$ws.Cells.Item(95, 14).Value = 4300                       # N95 - Precio mínimo
$ws.Cells.Item(95, 15).Value = 4500                       # O95 - Precio máximo
$ws.Cells.Item(95, 16).Value = 4397                       # P95 - Precio promedio ponderado
$ws.Cells.Item(95, 17).Value = "$/bandeja 2 kilos"       # Q95 - Unidad de comercialización
$ws.Cells.Item(95, 18).Value = "Provincia de Linares"    # R95 - Origen
$ws.Cells.Item(95, 19).Value = 2198                       # S95 - Precio $/Kg
$ws.Cells.Item(95, 20).Value = 2                          # T95 - Kg / unidad
